$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.497.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").Value = "'1.857.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'245.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").Value = "'0.6963"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.19%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "'0.07701"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").Value = "'23.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("D11").Value = "'0.07796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").Value = "'5.162"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").Value = "'1.851.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("E14").Value = "  +2.02%  "

$ws.Range("D15").Value = "'91.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "'6.355"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").Value = "'29.483.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("D18").Value = "'0.000008316"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").Value = "'2.100.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").Value = "'238.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("D21").Value = "'12.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +2.26%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").Value = "'160.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "'8.903"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("D28").Value = "'18.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'1.536"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'4.253"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("D31").Value = "'4.151"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "'1.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("D33").Value = "'0.05111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").Value = "'0.7774"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("D38").Value = "'1.317.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.26%  "

$ws.Range("D39").Value = "'0.01874"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("D41").Value = "'0.9532"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").Value = "'105.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("D43").Value = "'5.770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").Value = "'9.820"
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("D47").Value = "'2.000.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("D48").Value = "'0.5233"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "

$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").Value = "'63.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").Value = "'6.969"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
